$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 77.14286
$ws.Range("J5").Value = 50
$ws.Range("L5").Value = 50
$ws.Range("N5").Value = -280
$ws.Range("H17").Value = 982.3077
$ws.Range("J17").Value = 982.3077
$ws.Range("L17").Value = 2946.9231
$ws.Range("N17").Value = -3282.9231
$ws.Range("H33").Value = 2703.5454
$ws.Range("I33").Value = 2984.875
$ws.Range("J33").Value = 1953.3334
$ws.Range("K33").Value = 2984.875
$ws.Range("L33").Value = 1953.3334
$ws.Range("M33").Value = -2755.875
$ws.Range("N33").Value = -2411.3334
$ws.Range("H42").Value = 311.2
$ws.Range("I42").Value = 311.2
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 933.5999999999999
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -703.5999999999999
$ws.Range("H62").Value = 2160.6
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2160.6
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H100").Value = 1899
$ws.Range("J100").Value = 2403
$ws.Range("L100").Value = 2403
$ws.Range("N100").Value = -3485
$ws.Range("H116").Value = 7667.273
$ws.Range("I116").Value = 6161
$ws.Range("J116").Value = 8922.5
$ws.Range("K116").Value = 6161
$ws.Range("L116").Value = 8922.5
$ws.Range("M116").Value = -2719
$ws.Range("N116").Value = -15806.5
$ws.Range("H137").Value = 4746.737
$ws.Range("J137").Value = 8148.857
$ws.Range("L137").Value = 24446.571
$ws.Range("N137").Value = -29546.571
$ws.Range("H138").Value = 2139.3953
$ws.Range("J138").Value = 2871.84
$ws.Range("L138").Value = 8615.52
$ws.Range("N138").Value = -18895.52
$ws.Range("H141").Value = 3174.1333
$ws.Range("J141").Value = 2997.5
$ws.Range("L141").Value = 8992.5
$ws.Range("N141").Value = -19352.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12509170
$ws.Range("I74").Value = 19232232
$ws.Range("J74").Value = 23483.143
$ws.Range("K74").Value = 19232232
$ws.Range("L74").Value = 23483.143
$ws.Range("M74").Value = -19231358
$ws.Range("N74").Value = -25231.143
$ws.Range("H77").Value = 12509170
$ws.Range("I77").Value = 19232232
$ws.Range("J77").Value = 23483.143
$ws.Range("K77").Value = 96161160
$ws.Range("L77").Value = 117415.715
$ws.Range("M77").Value = -96156792
$ws.Range("N77").Value = -126151.715

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2724.7932
$ws.Range("J20").Value = 2380.7144
$ws.Range("L20").Value = 2380.7144
$ws.Range("N20").Value = -2874.7144
$ws.Range("H132").Value = 90714.28999999999
$ws.Range("J132").Value = 67500
$ws.Range("L132").Value = 67500
$ws.Range("N132").Value = -77620

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 620686
$ws.Range("I31").Value = 8413.387000000001
$ws.Range("K31").Value = 8413.387000000001
$ws.Range("M31").Value = -8118.387000000001
$ws.Range("H34").Value = 620686
$ws.Range("I34").Value = 8413.387000000001
$ws.Range("K34").Value = 8413.387000000001
$ws.Range("M34").Value = -8211.387000000001
$ws.Range("H99").Value = 3839
$ws.Range("I99").Value = 3902
$ws.Range("J99").Value = 3650
$ws.Range("K99").Value = 3902
$ws.Range("L99").Value = 3650
$ws.Range("M99").Value = -2404
$ws.Range("N99").Value = -6646
$ws.Range("H125").Value = 44500
$ws.Range("J125").Value = 44500
$ws.Range("L125").Value = 44500
$ws.Range("N125").Value = -49420
$ws.Range("H126").Value = 3839
$ws.Range("I126").Value = 3902
$ws.Range("J126").Value = 3650
$ws.Range("K126").Value = 11706
$ws.Range("L126").Value = 10950
$ws.Range("M126").Value = -9236
$ws.Range("N126").Value = -15890
$ws.Range("H132").Value = 1413.2307
$ws.Range("I132").Value = 1447.6666
$ws.Range("K132").Value = 4342.9998
$ws.Range("M132").Value = -1812.9998
$ws.Range("H134").Value = 1007116.2
$ws.Range("I134").Value = 1667689.6
$ws.Range("J134").Value = 16256
$ws.Range("K134").Value = 5003068.800000001
$ws.Range("L134").Value = 48768
$ws.Range("M134").Value = -5000533.800000001
$ws.Range("N134").Value = -53838

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 54328.668
$ws.Range("I44").Value = 54328.668
$ws.Range("K44").Value = 162986.004
$ws.Range("M44").Value = -162588.004
$ws.Range("H46").Value = 2019.8
$ws.Range("J46").Value = 2499.75
$ws.Range("L46").Value = 7499.25
$ws.Range("N46").Value = -7681.25
$ws.Range("H62").Value = 21666.666
$ws.Range("I62").Value = 21666.666
$ws.Range("K62").Value = 64999.99800000001
$ws.Range("M62").Value = -64313.99800000001
$ws.Range("H65").Value = 21666.666
$ws.Range("I65").Value = 21666.666
$ws.Range("K65").Value = 194999.994
$ws.Range("M65").Value = -191567.994
$ws.Range("H92").Value = 2502725
$ws.Range("J92").Value = 450
$ws.Range("L92").Value = 1350
$ws.Range("N92").Value = -3846
$ws.Range("H109").Value = 826.44446
$ws.Range("I109").Value = 826.44446
$ws.Range("K109").Value = 2479.33338
$ws.Range("M109").Value = -1439.33338
$ws.Range("H129").Value = 25718054
$ws.Range("I129").Value = 4514.25
$ws.Range("J129").Value = 66859716
$ws.Range("K129").Value = 13542.75
$ws.Range("L129").Value = 200579148
$ws.Range("M129").Value = -8542.75
$ws.Range("N129").Value = -200589148
$ws.Range("H140").Value = 127431.375
$ws.Range("I140").Value = 138708.95
$ws.Range("K140").Value = 416126.85
$ws.Range("M140").Value = -410946.85

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0
$ws.Range("H24").Value = 24829.625
$ws.Range("I24").Value = 2875
$ws.Range("J24").Value = 32147.834
$ws.Range("K24").Value = 2875
$ws.Range("L24").Value = 32147.834
$ws.Range("M24").Value = -2702
$ws.Range("N24").Value = -32493.834
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2317.2
$ws.Range("I22").Value = 2317.2
$ws.Range("K22").Value = 2317.2
$ws.Range("M22").Value = -2022.2
$ws.Range("H27").Value = 2317.2
$ws.Range("I27").Value = 2317.2
$ws.Range("K27").Value = 2317.2
$ws.Range("M27").Value = -2210.2
$ws.Range("H40").Value = 3881.5217
$ws.Range("I40").Value = 3339.4707
$ws.Range("J40").Value = 5417.3335
$ws.Range("K40").Value = 3339.4707
$ws.Range("L40").Value = 5417.3335
$ws.Range("M40").Value = -3203.4707
$ws.Range("N40").Value = -5689.3335
$ws.Range("H82").Value = 760.0769
$ws.Range("I82").Value = 542.6667
$ws.Range("J82").Value = 1249.25
$ws.Range("K82").Value = 542.6667
$ws.Range("L82").Value = 1249.25
$ws.Range("M82").Value = -181.6667
$ws.Range("N82").Value = -1971.25
$ws.Range("H85").Value = 760.0769
$ws.Range("I85").Value = 542.6667
$ws.Range("J85").Value = 1249.25
$ws.Range("K85").Value = 542.6667
$ws.Range("L85").Value = 1249.25
$ws.Range("M85").Value = 705.3333
$ws.Range("N85").Value = -3745.25
$ws.Range("H122").Value = 5138.522
$ws.Range("I122").Value = 4778.7334
$ws.Range("J122").Value = 5813.125
$ws.Range("K122").Value = 14336.2002
$ws.Range("L122").Value = 17439.375
$ws.Range("M122").Value = -11886.2002
$ws.Range("N122").Value = -22339.375
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H136").Value = 32786.977
$ws.Range("I136").Value = 5172.6665
$ws.Range("J136").Value = 82492.734
$ws.Range("K136").Value = 15517.9995
$ws.Range("L136").Value = 247478.202
$ws.Range("M136").Value = -12967.9995
$ws.Range("N136").Value = -252578.202

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = 0
$ws.Range("H132").Value = 1533.6061
$ws.Range("I132").Value = 1374.2593
$ws.Range("K132").Value = 4122.7779
$ws.Range("M132").Value = -1592.7779
$ws.Range("H136").Value = 1053.4117
$ws.Range("I136").Value = 1064.9286
$ws.Range("K136").Value = 3194.7858
$ws.Range("M136").Value = -644.7857999999997
